$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap country-name pairs in column A (shared-string table reorder in the
# source diff, expressed here simply as swapping the two cell values) ---

# Dominica (row 202) <-> Fiyi (row 203)
$ws.Cells.Item(202, 1).Value = "Fiyi"
$ws.Cells.Item(203, 1).Value = "Dominica"

# Islas Malvinas (row 207) <-> Groenlandia (row 208)
$ws.Cells.Item(207, 1).Value = "Groenlandia"
$ws.Cells.Item(208, 1).Value = "Islas Malvinas"

# Islas Virgenes Britanicas (row 213) <-> Papua Nueva Guinea (row 214)
$ws.Cells.Item(213, 1).Value = "Papua Nueva Guinea"
$ws.Cells.Item(214, 1).Value = "Islas Virgenes Britanicas"

# --- Updated "last refreshed" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 22 de Junio de 2020 a las 09:11"

# --- Refreshed per-country case numbers ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 2356715
$ws.Cells.Item(4, 3).Value = 58
$ws.Cells.Item(4, 4).Value = 980367
$ws.Cells.Item(4, 5).Value = 1254100
$ws.Cells.Item(4, 7).Value = 1
$ws.Cells.Item(4, 8).Value = 122248

# Row 38: Ucrania
$ws.Cells.Item(38, 2).Value = 37241
$ws.Cells.Item(38, 3).Value = 681
$ws.Cells.Item(38, 4).Value = 16642
$ws.Cells.Item(38, 5).Value = 19587
$ws.Cells.Item(38, 7).Value = 10
$ws.Cells.Item(38, 8).Value = 1012

# Row 52: Armenia
$ws.Cells.Item(52, 2).Value = 20588
$ws.Cells.Item(52, 3).Value = 320
$ws.Cells.Item(52, 4).Value = 9131
$ws.Cells.Item(52, 5).Value = 11097
$ws.Cells.Item(52, 7).Value = 10
$ws.Cells.Item(52, 8).Value = 360

# Row 67: Chequia
$ws.Cells.Item(67, 4).Value = 7505
$ws.Cells.Item(67, 5).Value = 2657

# Row 89: Hungria
$ws.Cells.Item(89, 2).Value = 4102
$ws.Cells.Item(89, 3).Value = 8
$ws.Cells.Item(89, 4).Value = 2590
$ws.Cells.Item(89, 5).Value = 940
$ws.Cells.Item(89, 7).Value = 2
$ws.Cells.Item(89, 8).Value = 572

# Row 130: Georgia
$ws.Cells.Item(130, 2).Value = 908
$ws.Cells.Item(130, 3).Value = 2
$ws.Cells.Item(130, 4).Value = 761
$ws.Cells.Item(130, 5).Value = 133

# Row 213: now Papua Nueva Guinea
$ws.Cells.Item(213, 4).Value = 8
$ws.Cells.Item(213, 8).Value = 0

# Row 214: now Islas Virgenes Britanicas
$ws.Cells.Item(214, 4).Value = 7
$ws.Cells.Item(214, 8).Value = 1
